# Updates cryptos list: refresh Price/Volume(1h) figures and reorder two
# coin pairs (InternetComputer/Filecoin and Kaspa/FraxShare swapped ranks).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.556.08"
$ws.Range("E2").Value = "  +3.18%  "
$ws.Range("D3").Value = "1.856.10"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("D5").Value = "'230.08"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").Value = "'0.610"
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("D7").Value = "'1.01"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'42.11"
$ws.Range("E8").Value = "  +10.53%  "
$ws.Range("D9").Value = "'0.307"
$ws.Range("E9").Value = "  +6.80%  "
$ws.Range("D10").Value = "'0.0692"
$ws.Range("E10").Value = "  +3.08%  "
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("D12").Value = "2.126.90"
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("D13").Value = "'11.48"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("D14").Value = "1.840.82"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "'0.673"
$ws.Range("E15").Value = "  +7.10%  "
$ws.Range("D16").Value = "'4.70"
$ws.Range("E16").Value = "  +6.63%  "
$ws.Range("D17").Value = "35.548.98"
$ws.Range("E17").Value = "  +3.17%  "
$ws.Range("D18").Value = "'70.05"
$ws.Range("E18").Value = "  +3.08%  "
$ws.Range("D19").Value = "'247.35"
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("D20").Value = "0.0₃0800"
$ws.Range("E20").Value = "  +4.06%  "
$ws.Range("D21").Value = "'12.15"
$ws.Range("E21").Value = "  +9.42%  "
$ws.Range("D22").Value = "'4.65"
$ws.Range("E22").Value = "  +13.69%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "'169.06"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").Value = "'7.93"
$ws.Range("E26").Value = "  +3.07%  "
$ws.Range("D27").Value = "'17.78"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").Value = "'0.122"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("D29").Value = "'1.39"
$ws.Range("E29").Value = "  +13.45%  "
$ws.Range("D30").Value = "'1.01"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").Value = "3.299.49"
$ws.Range("E31").Value = "  +35.80%  "
$ws.Range("D32").Value = "'0.0543"
$ws.Range("E32").Value = "  +5.77%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.93"
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.06"
$ws.Range("E34").Value = "  +5.92%  "
$ws.Range("D35").Value = "'1.88"
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("D36").Value = "'99.53"
$ws.Range("E36").Value = "  +22.11%  "
$ws.Range("D37").Value = "'0.698"
$ws.Range("E37").Value = "  +9.21%  "
$ws.Range("D38").Value = "1.365.14"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "'2.47"
$ws.Range("E39").Value = "  +6.41%  "
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("D41").Value = "'0.0194"
$ws.Range("E41").Value = "  +3.64%  "
$ws.Range("E42").Value = "  +6.95%  "
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("D44").Value = "'14.69"
$ws.Range("E44").Value = "  +7.80%  "
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").Value = "'2.81"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").Value = "'0.0521"
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'6.22"
$ws.Range("E48").Value = "  +8.34%  "
$ws.Range("D49").Value = "2.024.26"
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("D50").Value = "'1.01"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "'103.80"
$ws.Range("E51").Value = "  +1.76%  "
